# neue samples in die listen eingefuegt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (thin borders) of the last existing data row (86)
# down onto the three new rows before writing their values, so the new
# rows pick up the same visual style as the preceding entries.
$ws.Range("A86").Copy()
$ws.Range("A87:A89").PasteSpecial(-4122)
$ws.Range("B86").Copy()
$ws.Range("B87:B89").PasteSpecial(-4122)

# New sample rows appended to the table.
$ws.Range("A87").Value = 86
$ws.Range("B87").Value = "kickle"

$ws.Range("A88").Value = 87
$ws.Range("B88").Value = "SEA_subkick"

$ws.Range("A89").Value = 88
$ws.Range("B89").Value = "SEA_subkick2"

# Update the visible selection to match where the user ended up after
# adding the rows.
$ws.Range("D96").Select() | Out-Null
